$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update column F values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 795
$ws1.Range("F3").Value = 971
$ws1.Range("F4").Value = 755
$ws1.Range("F5").Value = 859
$ws1.Range("F6").Value = 422
$ws1.Range("F7").Value = 654
$ws1.Range("F10").Value = 671
$ws1.Range("F12").Value = 530
$ws1.Range("F15").Value = 737
$ws1.Range("F17").Value = 381
$ws1.Range("F21").Value = 122
$ws1.Range("F22").Value = 611
$ws1.Range("F24").Value = 876

# Sheet "演出" (Performances) - update column F values
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 333
$ws2.Range("F5").Value = 26
$ws2.Range("F8").Value = 234

# Sheet "全部类型" (All types) - update column F values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 333
$ws4.Range("F5").Value = 795
$ws4.Range("F6").Value = 971
$ws4.Range("F7").Value = 755
$ws4.Range("F8").Value = 859
$ws4.Range("F9").Value = 422
$ws4.Range("F10").Value = 654
$ws4.Range("F13").Value = 671
$ws4.Range("F15").Value = 26
$ws4.Range("F17").Value = 530
$ws4.Range("F21").Value = 737
$ws4.Range("F24").Value = 381
$ws4.Range("F27").Value = 234
$ws4.Range("F34").Value = 122
$ws4.Range("F35").Value = 611
$ws4.Range("F37").Value = 876
